# Update res_line/pl_mw.xlsx values for the "case with 380 kV done" run.
# Applies new computed values to Sheet1 (columns B-N, rows 2-25),
# matching the diff between the previous run and the 380 kV case.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (index A2 = 0)
$ws.Range("B2").Value = 0.5258102763072543
$ws.Range("C2").Value = 0.05157610560628711
$ws.Range("D2").Value = 0.07795294679323206
$ws.Range("E2").Value = 0.1262542582653552
$ws.Range("G2").Value = 1.322609435993741
$ws.Range("H2").Value = 1.209598308510806
$ws.Range("K2").Value = 0.5003313440759314
$ws.Range("M2").Value = 0.2801782571528904
$ws.Range("N2").Value = 2.216355202512304
# Row 3 (index A3 = 1)
$ws.Range("B3").Value = 0.4842161353644201
$ws.Range("C3").Value = 0.04496850578706812
$ws.Range("D3").Value = 0.07078566923591723
$ws.Range("E3").Value = 0.116341087167747
$ws.Range("G3").Value = 1.302129554721844
$ws.Range("H3").Value = 1.205095884187273
$ws.Range("K3").Value = 0.4547432370559079
$ws.Range("M3").Value = 0.2562795544629779
$ws.Range("N3").Value = 2.225051549865086
# Row 4 (index A4 = 2)
$ws.Range("B4").Value = 0.4589542006571321
$ws.Range("C4").Value = 0.04091563911713081
$ws.Range("D4").Value = 0.06642250110409975
$ws.Range("E4").Value = 0.1103314184761999
$ws.Range("G4").Value = 1.29027882387264
$ws.Range("H4").Value = 1.202855476744162
$ws.Range("K4").Value = 0.4269901258888638
$ws.Range("M4").Value = 0.2417577190666336
$ws.Range("N4").Value = 2.231035679859957
# Row 5 (index A5 = 3)
$ws.Range("B5").Value = 0.4487292503866911
$ws.Range("C5").Value = 0.03926504994738877
$ws.Range("D5").Value = 0.06465385431555148
$ws.Range("E5").Value = 0.1079015473853246
$ws.Range("G5").Value = 1.285631066629549
$ws.Range("H5").Value = 1.202074153235003
$ws.Range("K5").Value = 0.4157401353682815
$ws.Range("M5").Value = 0.2358778382335416
$ws.Range("N5").Value = 2.233636060642169
# Row 6 (index A6 = 4)
$ws.Range("B6").Value = 0.4470355975948337
$ws.Range("C6").Value = 0.03899102828306411
$ws.Range("D6").Value = 0.06436073641739881
$ws.Range("E6").Value = 0.1074992162586952
$ws.Range("G6").Value = 1.284870257167157
$ws.Range("H6").Value = 1.201952363367383
$ws.Range("K6").Value = 0.4138756768630003
$ws.Range("M6").Value = 0.2349037664237272
$ws.Range("N6").Value = 2.234077616182624
# Row 7 (index A7 = 5)
$ws.Range("B7").Value = 0.4588160222475608
$ws.Range("C7").Value = 0.04089337479774713
$ws.Range("D7").Value = 0.06639861062426178
$ws.Range("E7").Value = 0.1102985713184381
$ws.Range("G7").Value = 1.29021540838248
$ws.Range("H7").Value = 1.202844406630078
$ws.Range("K7").Value = 0.4268381634945513
$ws.Range("M7").Value = 0.2416782680052521
$ws.Range("N7").Value = 2.231070094763766
# Row 8 (index A8 = 6)
$ws.Range("B8").Value = 0.5114110004659551
$ws.Range("C8").Value = 0.0492968802525553
$ws.Range("D8").Value = 0.07547382901914546
$ws.Range("E8").Value = 0.1228200438766649
$ws.Range("G8").Value = 1.315397369169887
$ws.Range("H8").Value = 1.20793701369459
$ws.Range("K8").Value = 0.4845629584565643
$ws.Range("M8").Value = 0.2719062027331915
$ws.Range("N8").Value = 2.219219762286656
# Row 9 (index A9 = 7)
$ws.Range("B9").Value = 0.6167614052727401
$ws.Range("C9").Value = 0.06581390327141889
$ws.Range("D9").Value = 0.09357240399332056
$ws.Range("E9").Value = 0.1479990570136209
$ws.Range("G9").Value = 1.370552455800862
$ws.Range("H9").Value = 1.222089820241564
$ws.Range("K9").Value = 0.5996690746065099
$ws.Range("M9").Value = 0.3324088312989204
$ws.Range("N9").Value = 2.201108333057959
# Row 10 (index A10 = 8)
$ws.Range("B10").Value = 0.6955382876831777
$ws.Range("C10").Value = 0.07797903946826068
$ws.Range("D10").Value = 0.1070605238471671
$ws.Range("E10").Value = 0.166898689211834
$ws.Range("G10").Value = 1.414640326666984
$ws.Range("H10").Value = 1.235041228300702
$ws.Range("K10").Value = 0.6854350325540963
$ws.Range("M10").Value = 0.3776386000478595
$ws.Range("N10").Value = 2.190946084181192
# Row 11 (index A11 = 9)
$ws.Range("B11").Value = 0.7316812235734744
$ws.Range("C11").Value = 0.08352140494045557
$ws.Range("D11").Value = 0.1132396932661379
$ws.Range("E11").Value = 0.1755880424070071
$ws.Range("G11").Value = 1.435481569639023
$ws.Range("H11").Value = 1.241490739572896
$ws.Range("K11").Value = 0.7247202267155046
$ws.Range("M11").Value = 0.3983908861920611
$ws.Range("N11").Value = 2.187009901303057
# Row 12 (index A12 = 10)
$ws.Range("B12").Value = 0.7454120974611271
$ws.Range("C12").Value = 0.08562146457222752
$ws.Range("D12").Value = 0.1155859161881665
$ws.Range("E12").Value = 0.1788919982743167
$ws.Range("G12").Value = 1.443487276549575
$ws.Range("H12").Value = 1.244013430161232
$ws.Range("K12").Value = 0.7396357655707959
$ws.Range("M12").Value = 0.4062751686865482
$ws.Range("N12").Value = 2.185618457278139
# Row 13 (index A13 = 11)
$ws.Range("B13").Value = 0.7424529334719239
$ws.Range("C13").Value = 0.08516912066431814
$ws.Range("D13").Value = 0.1150803335245598
$ws.Range("E13").Value = 0.1781798279134463
$ws.Range("G13").Value = 1.441758042638554
$ws.Range("H13").Value = 1.243466544818062
$ws.Range("K13").Value = 0.7364216943028907
$ws.Range("M13").Value = 0.4045759918866239
$ws.Range("N13").Value = 2.18591371622297
# Row 14 (index A14 = 12)
$ws.Range("B14").Value = 0.732809981309515
$ws.Range("C14").Value = 0.083694151870219
$ws.Range("D14").Value = 0.1134325916563199
$ws.Range("E14").Value = 0.1758595889249364
$ws.Range("G14").Value = 1.436137923989122
$ws.Range("H14").Value = 1.241696670439438
$ws.Range("K14").Value = 0.7259465520180868
$ws.Range("M14").Value = 0.3990390107276269
$ws.Range("N14").Value = 2.186893438115291
# Row 15 (index A15 = 13)
$ws.Range("B15").Value = 0.7269091733141124
$ws.Range("C15").Value = 0.08279086079059539
$ws.Range("D15").Value = 0.1124241254733107
$ws.Range("E15").Value = 0.1744401403293807
$ws.Range("G15").Value = 1.432710249493994
$ws.Range("H15").Value = 1.240623047293496
$ws.Range("K15").Value = 0.7195353239934263
$ws.Range("M15").Value = 0.3956508267718988
$ws.Range("N15").Value = 2.187506462965445
# Row 16 (index A16 = 14)
$ws.Range("B16").Value = 0.6931824544190874
$ws.Range("C16").Value = 0.07761700841743391
$ws.Range("D16").Value = 0.1066575790872832
$ws.Range("E16").Value = 0.1663326912833298
$ws.Range("G16").Value = 1.413294161711917
$ws.Range("H16").Value = 1.234630978456636
$ws.Range("K16").Value = 0.6828731068165723
$ws.Range("M16").Value = 0.3762859845555013
$ws.Range("N16").Value = 2.191217171304174
# Row 17 (index A17 = 15)
$ws.Range("B17").Value = 0.6725709532546489
$ws.Range("C17").Value = 0.07444522189258862
$ws.Range("D17").Value = 0.1031311445434113
$ws.Range("E17").Value = 0.1613827509842238
$ws.Range("G17").Value = 1.401584642175635
$ws.Range("H17").Value = 1.231098045836319
$ws.Range("K17").Value = 0.6604513365452647
$ws.Range("M17").Value = 0.3644518968540922
$ws.Range("N17").Value = 2.193669699223705
# Row 18 (index A18 = 16)
$ws.Range("B18").Value = 0.6607446415588925
$ws.Range("C18").Value = 0.07262167844957901
$ws.Range("D18").Value = 0.1011069092662922
$ws.Range("E18").Value = 0.158544301292892
$ws.Range("G18").Value = 1.394923528151224
$ws.Range("H18").Value = 1.229118492816497
$ws.Range("K18").Value = 0.6475803074610553
$ws.Range("M18").Value = 0.3576618830092144
$ws.Range("N18").Value = 2.19514493964347
# Row 19 (index A19 = 17)
$ws.Range("B19").Value = 0.6567454098041594
$ws.Range("C19").Value = 0.07200439023364424
$ws.Range("D19").Value = 0.100422235986386
$ws.Range("E19").Value = 0.1575847237604364
$ws.Range("G19").Value = 1.392680861200915
$ws.Range("H19").Value = 1.228457260998113
$ws.Range("K19").Value = 0.6432267485680541
$ws.Range("M19").Value = 0.3553657473443792
$ws.Range("N19").Value = 2.195655516872975
# Row 20 (index A20 = 18)
$ws.Range("B20").Value = 0.6747620925435172
$ws.Range("C20").Value = 0.07478278223510415
$ws.Range("D20").Value = 0.1035061170960745
$ws.Range("E20").Value = 0.1619087858767259
$ws.Range("G20").Value = 1.402823488407847
$ws.Range("H20").Value = 1.231468698048019
$ws.Range("K20").Value = 0.6628355434584989
$ws.Range("M20").Value = 0.3657099307054921
$ws.Range("N20").Value = 2.193401933812865
# Row 21 (index A21 = 19)
$ws.Range("B21").Value = 0.7356411454967429
$ws.Range("C21").Value = 0.08412735040815278
$ws.Range("D21").Value = 0.1139164016852021
$ws.Range("E21").Value = 0.1765407310717393
$ws.Range("G21").Value = 1.437785601444062
$ws.Range("H21").Value = 1.242214341922988
$ws.Range("K21").Value = 0.7290222920625808
$ws.Range("M21").Value = 0.4006646518823729
$ws.Range("N21").Value = 2.186602977525141
# Row 22 (index A22 = 20)
$ws.Range("B22").Value = 0.7756875985171519
$ws.Range("C22").Value = 0.0902421003896734
$ws.Range("D22").Value = 0.1207569114959313
$ws.Range("E22").Value = 0.1861822954207497
$ws.Range("G22").Value = 1.461297696294537
$ws.Range("H22").Value = 1.249705965270948
$ws.Range("K22").Value = 0.7725071329994648
$ws.Range("M22").Value = 0.4236603916699693
$ws.Range("N22").Value = 2.182737285317131
# Row 23 (index A23 = 21)
$ws.Range("B23").Value = 0.7542903380848998
$ws.Range("C23").Value = 0.08697782698899914
$ws.Range("D23").Value = 0.1171026130422916
$ws.Range("E23").Value = 0.1810291118592247
$ws.Range("G23").Value = 1.448688034178389
$ws.Range("H23").Value = 1.245664595371693
$ws.Range("K23").Value = 0.7492774989346174
$ws.Range("M23").Value = 0.4113732055137547
$ws.Range("N23").Value = 2.184747487339919
# Row 24 (index A24 = 22)
$ws.Range("B24").Value = 0.6737714052423769
$ws.Range("C24").Value = 0.07463017134044492
$ws.Range("D24").Value = 0.1033365821738386
$ws.Range("E24").Value = 0.1616709426745544
$ws.Range("G24").Value = 1.402263185387795
$ws.Range("H24").Value = 1.231300965521825
$ws.Range("K24").Value = 0.6617575826594759
$ws.Range("M24").Value = 0.3651411313972233
$ws.Range("N24").Value = 2.193522787305838
# Row 25 (index A25 = 23)
$ws.Range("B25").Value = 0.5880217089342636
$ws.Range("C25").Value = 0.06134084134396289
$ws.Range("D25").Value = 0.08864320381172774
$ws.Range("E25").Value = 0.1411185635444241
$ws.Range("G25").Value = 1.355008751530363
$ws.Range("H25").Value = 1.217813706960129
$ws.Range("K25").Value = 0.5683220290439124
$ws.Range("M25").Value = 0.3159068272816086
$ws.Range("N25").Value = 2.205457161484006

Write-Host "Updated 216 cells on $($ws.Name)"
